$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Kupci": swap the sample customers for new placeholder data.
# (Edited first so the new shared strings land in the same order the
# original author typed them in: company names, then tax numbers.)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Kupci")
$ws2.Range("C2").Value = "Kompanija1 d.o.o."
$ws2.Range("C3").Value = "Kompanija2 d.d"
$ws2.Range("B2").Value = "21111111114"
$ws2.Range("B3").Value = "81111111110"
$ws2.Range("C2").Select()

# ------------------------------------------------------------------
# Sheet "Racuni": same customer tax-number swap for the invoice rows.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Racuni")
$ws3.Range("G2").Value = "21111111114"
$ws3.Range("G3").Value = "21111111114"
$ws3.Range("G4").Value = "81111111110"
$ws3.Range("C10").Select()

# ------------------------------------------------------------------
# Sheet "Zaglavlje": drop the "sastavio_*" (prepared-by) columns and
# add two new "opz_*" total columns.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Zaglavlje")

# The hyperlink engine only supports clearing ALL hyperlinks on a
# sheet at once, so drop them all and re-create the one that survives
# (column H, "email") after the column shuffle below.
$ws1.Range("A1").Hyperlinks.Delete()

# Remove columns I:M (sastavio_ime, sastavio_prezime, sastavio_tel,
# sastavio_fax, sastavio_email) - this shifts N/O (na_dan,
# nisu_naplaceni_do) left into I/J.
$ws1.Range("I1:M1").EntireColumn.Delete()

# Re-create the hyperlink for the "email" column (now still H).
$ws1.Hyperlinks.Add($ws1.Range("H2"), "hrvoje.jesenovic@gmail.com")

# Add the two new trailing columns with their header + default value.
$ws1.Range("K1").Value = "opz_ukupan_iznos_racuna_s_pdv"
$ws1.Range("L1").Value = "opz_ukupan_iznos_pdv"
$ws1.Range("K1:L1").Font.Bold = $true
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = 0
$ws1.Range("K2:L2").NumberFormat = "#,##0.00\ _H_R_K"

# ------------------------------------------------------------------
# View state: "Zaglavlje" becomes the active/selected sheet again.
# ------------------------------------------------------------------
$ws1.Select()
$ws1.Range("F15").Select()
